$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting rows 61-65 down to 62-66
$ws.Rows.Item(61).Insert()

# Row 61: A61 keeps the #NULL! error (style carried over from insert), B61 becomes "elz93"
$ws.Range("A61").Value = "#NULL!"
$ws.Range("B61").Value = "elz93"

